# Add a new "ValuesHistory" worksheet right after "Values", carrying a
# PreviousID-augmented view of the Values table (3 sample rows), matching
# the look & feel of the existing sheets (reusing their cell formatting).

$wb = $excel.ActiveWorkbook

$valuesSheet = $wb.Worksheets.Item("Values")
$keysSheet   = $wb.Worksheets.Item("Keys")
$storeSheet  = $wb.Worksheets.Item("Store")

# --- Create the new sheet right after "Values" -----------------------------
$newSheet = $wb.Worksheets.Add($null, $valuesSheet)
$newSheet.Name = "ValuesHistory"
$newSheet.Tab.Color = 5263615   # matches the other data sheets' red tab (FFFF5050)

# --- Header row (reuse styles from the Values / Keys sheets) ---------------
$valuesSheet.Range("A1:D1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$valuesSheet.Range("F1:I1").Copy()
$newSheet.Range("F1").PasteSpecial(-4122)
$keysSheet.Range("D1").Copy()
$newSheet.Range("E1").PasteSpecial(-4122)

# --- Data-row styles (rows 2-4), reuse Values row 2 formatting -------------
$valuesSheet.Range("A2:D2").Copy()
$newSheet.Range("A2:D4").PasteSpecial(-4122)
$valuesSheet.Range("F2:I2").Copy()
$newSheet.Range("F2:I4").PasteSpecial(-4122)

# --- Column widths (best effort match of the source widths) ----------------
$newSheet.Columns("A:E").ColumnWidth = 12.86
$newSheet.Columns("F").ColumnWidth = 19.63
$newSheet.Columns("G:H").ColumnWidth = 19.85
$newSheet.Columns("I").ColumnWidth = 12.86

# --- Header values -----------------------------------------------------
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "KeyID"
$newSheet.Range("C1").Value = "FieldID"
$newSheet.Range("D1").Value = "Value"
$newSheet.Range("E1").Value = "PreviousID"
$newSheet.Range("F1").Value = "CreationTime"
$newSheet.Range("G1").Value = "ModificationTime"
$newSheet.Range("H1").Value = "DeletionTime"
$newSheet.Range("I1").Value = "CommitID"

# --- Row 2 -------------------------------------------------------------
$newSheet.Range("A2").Value = "5A86F601F791F37D70A1A3B8FEDEB7BF2C7B88D4"
$newSheet.Range("B2").Value = "00BF313399C32EE1563AC7BD598236C359126679"
$newSheet.Range("C2").Value = "4C90630588DA709A3007B7EE0FB7DDFD9159BE90"
$newSheet.Range("D2").Value = 70812
$newSheet.Range("E2").Value = "5A86F601F791F37D70A1A3B8FEDEB7BF2C7B88D4"
$newSheet.Range("F2").Value = 46012.915267939818
$newSheet.Range("G2").Value = 46012.915267939818
$newSheet.Range("H2").Value = 0
$newSheet.Range("I2").Value = "135EC372181380110551F971B7FA0C4703A3739F"

# --- Row 3 -------------------------------------------------------------
$newSheet.Range("A3").Value = "3028B7A9A3B110C6AA6640127F96D3310E5137D4"
$newSheet.Range("B3").Value = "00D3E1214778099D4FA23326203B06AC6E46DC43"
$newSheet.Range("C3").Value = "4C90630588DA709A3007B7EE0FB7DDFD9159BE90"
$newSheet.Range("D3").Value = 95334
$newSheet.Range("E3").Value = "3028B7A9A3B110C6AA6640127F96D3310E5137D4"
$newSheet.Range("F3").Value = 46012.915267824072
$newSheet.Range("G3").Value = 46012.915267824072
$newSheet.Range("H3").Value = 0
$newSheet.Range("I3").Value = "135EC372181380110551F971B7FA0C4703A3739F"

# --- Row 4 -------------------------------------------------------------
$newSheet.Range("A4").Value = "9598AED755869D9DE18753D2B2B0966604B1B825"
$newSheet.Range("B4").Value = "00DE27E4B7421EC9E52AD33D82AFA4855B8DC64D"
$newSheet.Range("C4").Value = "4C90630588DA709A3007B7EE0FB7DDFD9159BE90"
$newSheet.Range("D4").Value = 54771
$newSheet.Range("E4").Value = "9598AED755869D9DE18753D2B2B0966604B1B825"
$newSheet.Range("F4").Value = 46012.915268171295
$newSheet.Range("G4").Value = 46012.915268171295
$newSheet.Range("H4").Value = 0
$newSheet.Range("I4").Value = "135EC372181380110551F971B7FA0C4703A3739F"

# --- Restore "Store" as the active/selected sheet ---------------------
$storeSheet.Activate()
